# Collection_methods.xlsx maintenance edit:
#   - Sheet1 -> CollectionMethods (the workbook's _FilterDatabase defined
#     name tracks the sheet automatically once it is renamed)
#   - Row heights refreshed to match the new default font metrics; rows
#     that had an explicit (non custom-height) autofit height get a new
#     autofit height, and the two rows that now fit within a single
#     default-height line (1 and 33) go back to the sheet default.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "CollectionMethods"

# Rows that now collapse back onto the (new) default row height.
$ws.Rows.Item(1).AutoFit()
$ws.Rows.Item(33).AutoFit()

# Rows whose autofit height shrank from two lines @15pt (30) to two
# lines @14.5pt (29).
$rowsTo29 = @(2,3,4,11,12,15,17,18,22,30,31,32,34,35,36,37)
foreach ($r in $rowsTo29) {
    $ws.Rows.Item($r).RowHeight = 29
}

# Rows whose autofit height shrank from three lines @15pt (45) to three
# lines @14.5pt (43.5).
$rowsTo435 = @(7,8,9,10,13,14,25,26,27,28,29)
foreach ($r in $rowsTo435) {
    $ws.Rows.Item($r).RowHeight = 43.5
}
